$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Colors (OLE BGR) matching the workbook's existing fill palette
$colorPink       = 12695295  # FFB6C1 - "-" removed row
$colorGreen      = 12713910  # B6FFC1 - "+" added row
$colorPaleBlue   = 16771304  # E8E8FF
$colorPalePink   = 14867967  # FFDDE2 - Date/Origin column accent
$colorPaleGreen  = 14876637  # DDFFE2 - Shape column accent
$colorPaleYellow = 14548991  # FFFFDD - "c-"/"c+" data row accent

# xlPasteValues
$xlPasteValues = -4163

# Insert a new row above the current row 13, shifting the trailing rows down
$ws.Rows.Item(13).Insert()

# --- Row 12 becomes the "c-" (changed/old) variant of the ID=4 entry ---
$ws.Cells.Item(12, 1).Value = "c-"
$ws.Cells.Item(12, 1).Interior.Color = $colorPaleYellow

# "4" looks numeric, so write it via a text formula then flatten to a
# literal value with PasteSpecial - this keeps it stored as text instead
# of Excel auto-converting it to a number.
$ws.Cells.Item(12, 3).Formula = "=""4"""
$ws.Cells.Item(12, 3).Copy()
$ws.Cells.Item(12, 3).PasteSpecial($xlPasteValues)
$ws.Cells.Item(12, 3).Interior.Color = $colorPaleYellow

$ws.Cells.Item(12, 4).Value = "dried mango"
$ws.Cells.Item(12, 4).Interior.Color = $colorPaleYellow

$ws.Cells.Item(12, 7).Value = ""
$ws.Cells.Item(12, 7).Interior.Color = $colorPaleYellow

# --- Row 13 (new) becomes the "c+" (changed/new) variant of the ID=4 entry ---
$ws.Cells.Item(13, 1).Value = "c+"
$ws.Cells.Item(13, 1).Interior.Color = $colorPaleYellow

$ws.Cells.Item(13, 2).Value = "2022-11-21 00:00:00"
$ws.Cells.Item(13, 2).Interior.Color = $colorPalePink

$ws.Cells.Item(13, 3).Formula = "=""4"""
$ws.Cells.Item(13, 3).Copy()
$ws.Cells.Item(13, 3).PasteSpecial($xlPasteValues)
$ws.Cells.Item(13, 3).Interior.Color = $colorPaleYellow

$ws.Cells.Item(13, 4).Value = "dried mango"
$ws.Cells.Item(13, 4).Interior.Color = $colorPaleYellow

$ws.Cells.Item(13, 5).Value = "Thailand"
$ws.Cells.Item(13, 5).Interior.Color = $colorPalePink

$ws.Cells.Item(13, 6).Value = "flat"
$ws.Cells.Item(13, 6).Interior.Color = $colorPaleGreen

$ws.Cells.Item(13, 7).Value = "orange"
$ws.Cells.Item(13, 7).Interior.Color = $colorGreen
